$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion message in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$nl = [char]10
$newText = "Conversión del día 💰" + $nl + `
    "✅ Dólar paralelo: 68" + $nl + `
    "" + $nl + `
    "Binance" + $nl + `
    "✅ 1000 Bs = 9.64 = 39026.88 pesos" + $nl + `
    "✅ 39026.88 pesos = 9.6 = 940.95 Bs" + $nl + `
    "" + $nl + `
    "Promedio competencia" + $nl + `
    "✅ Tasa pesos: 20" + $nl + `
    "✅ Tasa Bs: 20" + $nl + `
    "✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 103.749
$wsTasas.Range("O10").Value = 4049
$wsTasas.Range("N12").Value = 4064.99
$wsTasas.Range("O12").Value = 98.008
